# Populate the new "Valeur Absolue" / "Angle" mini-table (rows 16-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Valeur Absolue"
$ws.Range("B16").Value = "Angle"

$ws.Range("A17").Value = 445
$ws.Range("B17").Value = 180

$ws.Range("A18").Value = 230
$ws.Range("B18").Value = 90

$ws.Range("A19").Value = 338
$ws.Range("B19").Value = 135

$ws.Range("A20").Value = 115
$ws.Range("B20").Value = 45

# Match the saved selection state of the sheet
[void]$ws.Range("E27").Select()
